$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "363÷9=" "304÷9="
Replace-Text "568÷3=" "173÷7="
Replace-Text "763÷7=" "779÷5="
Replace-Text "150÷4=" "382÷6="
Replace-Text "953÷3=" "882÷4="
Replace-Text "932÷8=" "947÷5="
Replace-Text "845÷4=" "391÷4="
Replace-Text "399÷5=" "809÷7="
Replace-Text "378÷6=" "530÷5="
Replace-Text "692÷8=" "558÷6="
Replace-Text "972÷4=" "106÷4="
Replace-Text "119÷8=" "643÷7="
Replace-Text "847÷7=" "805÷6="
Replace-Text "838÷8=" "732÷5="
Replace-Text "937÷7=" "876÷6="
Replace-Text "687÷2=" "421÷6="
Replace-Text "926÷8=" "414÷8="
Replace-Text "160÷3=" "681÷3="
Replace-Text "470÷3=" "649÷8="
Replace-Text "212÷5=" "127÷2="
Replace-Text "585÷2=" "855÷9="
Replace-Text "805÷3=" "746÷4="
Replace-Text "668÷3=" "658÷2="
Replace-Text "586÷4=" "416÷2="
Replace-Text "258÷4=" "453÷5="
